# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the H:N "market price" columns
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 793.8333
$ws.Range("I18").Value = 852.6
$ws.Range("K18").Value = 852.6
$ws.Range("M18").Value = -568.6

$ws.Range("H62").Value = 33336970
$ws.Range("I62").Value = 44447696
$ws.Range("K62").Value = 44447696
$ws.Range("M62").Value = -44447072

$ws.Range("H65").Value = 33336970
$ws.Range("I65").Value = 44447696
$ws.Range("K65").Value = 222238480
$ws.Range("M65").Value = -222235360

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2700.3635
$ws.Range("I32").Value = 2739.5813
$ws.Range("J32").Value = 1014
$ws.Range("K32").Value = 2739.5813
$ws.Range("L32").Value = 1014
$ws.Range("M32").Value = -2452.5813
$ws.Range("N32").Value = -1588

$ws.Range("H61").Value = 5186.56
$ws.Range("I61").Value = 3961.3157
$ws.Range("K61").Value = 3961.3157
$ws.Range("M61").Value = -3749.3157

$ws.Range("H74").Value = 3291.5862
$ws.Range("I74").Value = 1217.1177
$ws.Range("K74").Value = 1217.1177
$ws.Range("M74").Value = -343.1177

$ws.Range("H77").Value = 3291.5862
$ws.Range("I77").Value = 1217.1177
$ws.Range("K77").Value = 6085.5885
$ws.Range("M77").Value = -1717.5885

$ws.Range("H102").Value = 3460.875
$ws.Range("I102").Value = 3460.875
$ws.Range("K102").Value = 3460.875
$ws.Range("M102").Value = -1838.875

$ws.Range("H110").Value = 77060000
$ws.Range("I110").Value = 77060000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 77060000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -77057955
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 3794.4092
$ws.Range("I132").Value = 2804.1875
$ws.Range("K132").Value = 8412.5625
$ws.Range("M132").Value = -5882.5625

$ws.Range("H136").Value = 5186.56
$ws.Range("I136").Value = 3961.3157
$ws.Range("K136").Value = 11883.9471
$ws.Range("M136").Value = -9333.947100000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4018.5715
$ws.Range("J64").Value = 4189
$ws.Range("L64").Value = 4189
$ws.Range("N64").Value = -4639

$ws.Range("H67").Value = 4018.5715
$ws.Range("J67").Value = 4189
$ws.Range("L67").Value = 4189
$ws.Range("N67").Value = -5749

$ws.Range("H86").Value = 69617.10000000001
$ws.Range("J86").Value = 102403.65
$ws.Range("L86").Value = 102403.65
$ws.Range("N86").Value = -104649.65

$ws.Range("H89").Value = 69617.10000000001
$ws.Range("J89").Value = 102403.65
$ws.Range("L89").Value = 512018.25
$ws.Range("N89").Value = -523250.25

$ws.Range("H105").Value = 125002000
$ws.Range("I105").Value = 142859070
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 142859070
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -142857323
$ws.Range("N105").Value = -5994

$ws.Range("H134").Value = 3589.72
$ws.Range("I134").Value = 873.05554
$ws.Range("K134").Value = 2619.16662
$ws.Range("M134").Value = -84.16661999999997


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3925.2083
$ws.Range("J31").Value = 6389.5
$ws.Range("L31").Value = 6389.5
$ws.Range("N31").Value = -6979.5

$ws.Range("H34").Value = 3925.2083
$ws.Range("J34").Value = 6389.5
$ws.Range("L34").Value = 6389.5
$ws.Range("N34").Value = -6793.5

$ws.Range("H38").Value = 16721.428
$ws.Range("I38").Value = 12210
$ws.Range("K38").Value = 12210
$ws.Range("M38").Value = -11833

$ws.Range("H42").Value = 9487.5
$ws.Range("I42").Value = 9487.5
$ws.Range("K42").Value = 9487.5
$ws.Range("M42").Value = -8894.5

$ws.Range("H46").Value = 16721.428
$ws.Range("I46").Value = 12210
$ws.Range("K46").Value = 12210
$ws.Range("M46").Value = -11999

$ws.Range("H86").Value = 10220.625
$ws.Range("I86").Value = 8351.4
$ws.Range("J86").Value = 13336
$ws.Range("K86").Value = 8351.4
$ws.Range("L86").Value = 13336
$ws.Range("M86").Value = -7228.4
$ws.Range("N86").Value = -15582

$ws.Range("H89").Value = 10220.625
$ws.Range("I89").Value = 8351.4
$ws.Range("J89").Value = 13336
$ws.Range("K89").Value = 41757
$ws.Range("L89").Value = 66680
$ws.Range("M89").Value = -36141
$ws.Range("N89").Value = -77912

$ws.Range("H105").Value = 2743.3333
$ws.Range("I105").Value = 2740
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 2740
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = -993
$ws.Range("N105").Value = -6244

$ws.Range("H107").Value = 1911.8182
$ws.Range("I107").Value = 1538.3334
$ws.Range("K107").Value = 1538.3334
$ws.Range("M107").Value = 381.6666


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 126.0625
$ws.Range("I33").Value = 106.09091
$ws.Range("J33").Value = 170
$ws.Range("K33").Value = 636.5454599999999
$ws.Range("L33").Value = 1020
$ws.Range("M33").Value = -353.5454599999999
$ws.Range("N33").Value = -1586

$ws.Range("H34").Value = 2741.8572
$ws.Range("I34").Value = 491.66666
$ws.Range("K34").Value = 1474.99998
$ws.Range("M34").Value = -1390.99998

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H44").Value = 3837.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 3837.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 11512.5
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -12308.5

$ws.Range("H103").Value = 5369
$ws.Range("I103").Value = 5711.25
$ws.Range("J103").Value = 4000
$ws.Range("K103").Value = 17133.75
$ws.Range("L103").Value = 12000
$ws.Range("M103").Value = -16254.75
$ws.Range("N103").Value = -13758

$ws.Range("H137").Value = 2976.0833
$ws.Range("I137").Value = 1352.1666
$ws.Range("K137").Value = 4056.4998
$ws.Range("M137").Value = 1043.5002


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 27902.467
$ws.Range("J24").Value = 18048.818
$ws.Range("L24").Value = 18048.818
$ws.Range("N24").Value = -18394.818

$ws.Range("H80").Value = 2836.6667
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2836.6667
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 585.9583
$ws.Range("I97").Value = 589.45
$ws.Range("J97").Value = 568.5
$ws.Range("K97").Value = 589.45
$ws.Range("L97").Value = 568.5
$ws.Range("M97").Value = -93.45000000000005
$ws.Range("N97").Value = -1560.5

$ws.Range("H126").Value = 2731.1667
$ws.Range("I126").Value = 2337.0908
$ws.Range("K126").Value = 7011.2724
$ws.Range("M126").Value = -4541.2724


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 487
$ws.Range("I55").Value = 339.54544
$ws.Range("J55").Value = 1298
$ws.Range("K55").Value = 339.54544
$ws.Range("L55").Value = 1298
$ws.Range("M55").Value = -166.54544
$ws.Range("N55").Value = -1644

$ws.Range("H82").Value = 2037.5333
$ws.Range("I82").Value = 1458.875
$ws.Range("J82").Value = 2698.8572
$ws.Range("K82").Value = 1458.875
$ws.Range("L82").Value = 2698.8572
$ws.Range("M82").Value = -1097.875
$ws.Range("N82").Value = -3420.8572

$ws.Range("H85").Value = 2037.5333
$ws.Range("I85").Value = 1458.875
$ws.Range("J85").Value = 2698.8572
$ws.Range("K85").Value = 1458.875
$ws.Range("L85").Value = 2698.8572
$ws.Range("M85").Value = -210.875
$ws.Range("N85").Value = -5194.8572

$ws.Range("H122").Value = 5046.2666
$ws.Range("I122").Value = 4290.364
$ws.Range("J122").Value = 7125
$ws.Range("K122").Value = 12871.092
$ws.Range("L122").Value = 21375
$ws.Range("M122").Value = -10421.092
$ws.Range("N122").Value = -26275

$ws.Range("H132").Value = 5960.185
$ws.Range("I132").Value = 4774.2856
$ws.Range("K132").Value = 14322.8568
$ws.Range("M132").Value = -11792.8568


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17950.555
$ws.Range("I62").Value = 34763
$ws.Range("J62").Value = 4500.6
$ws.Range("K62").Value = 34763
$ws.Range("L62").Value = 4500.6
$ws.Range("M62").Value = -34139
$ws.Range("N62").Value = -5748.6

$ws.Range("H65").Value = 17950.555
$ws.Range("I65").Value = 34763
$ws.Range("J65").Value = 4500.6
$ws.Range("K65").Value = 173815
$ws.Range("L65").Value = 22503
$ws.Range("M65").Value = -170695
$ws.Range("N65").Value = -28743

$ws.Range("H99").Value = 54810.668
$ws.Range("I99").Value = 54810.668
$ws.Range("K99").Value = 54810.668
$ws.Range("M99").Value = -51815.668

$ws.Range("H107").Value = 4044.4443
$ws.Range("I107").Value = 4317.3335
$ws.Range("J107").Value = 3498.6667
$ws.Range("K107").Value = 12952.0005
$ws.Range("L107").Value = 10496.0001
$ws.Range("M107").Value = -11032.0005
$ws.Range("N107").Value = -14336.0001

$ws.Range("H113").Value = 1818.5454
$ws.Range("I113").Value = 930
$ws.Range("K113").Value = 2790
$ws.Range("M113").Value = -620

$ws.Range("H126").Value = 4152.9443
$ws.Range("I126").Value = 3109.5625
$ws.Range("K126").Value = 9328.6875
$ws.Range("M126").Value = -6858.6875


Write-Host "Applied all market-price cell updates."
